$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 212 (shifts old rows 212-316 down to 213-317)
$ws.Rows.Item(212).Insert()

# Populate the newly inserted row 212 with the new data record
$ws.Cells.Item(212, 1).Value = 3
$ws.Cells.Item(212, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(212, 3).Value = "Coquimbo"
$ws.Cells.Item(212, 4).Value = 44609
$ws.Cells.Item(212, 5).Value = 5
$ws.Cells.Item(212, 6).Value = 100112031
$ws.Cells.Item(212, 7).Value = "Poroto verde"
$ws.Cells.Item(212, 8).Value = "Magnum"
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 65
$ws.Cells.Item(212, 11).Value = 23000
$ws.Cells.Item(212, 12).Value = 25000
$ws.Cells.Item(212, 13).Value = 23923
$ws.Cells.Item(212, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(212, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(212, 16).Value = 957
$ws.Cells.Item(212, 17).Value = 25
$ws.Cells.Item(212, 18).Value = "Hortaliza"

Write-Output "done"
